$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 216 (E216, F216) ---
$ws.Range("E216").Value = 4.4478
$ws.Range("F216").Value = 4.4478

# --- Copy the formatting (date style) from A216 down to the new A217:A219 cells ---
$ws.Range("A216").Copy()
$ws.Range("A217:A219").PasteSpecial(-4122)  # xlPasteFormats

# --- Add new row 217 ---
$ws.Range("A217").Value = 45047.33333333334
$ws.Range("B217").Value = "FX_IDC:USDRON"
$ws.Range("C217").Value = 4.4295
$ws.Range("D217").Value = 4.66844
$ws.Range("E217").Value = 4.4169
$ws.Range("F217").Value = 4.6422
$ws.Range("G217").Value = 0

# --- Add new row 218 ---
$ws.Range("A218").Value = 45078.33333333334
$ws.Range("B218").Value = "FX_IDC:USDRON"
$ws.Range("C218").Value = 4.6421
$ws.Range("D218").Value = 4.65664
$ws.Range("E218").Value = 4.5035
$ws.Range("F218").Value = 4.5365
$ws.Range("G218").Value = 0

# --- Add new row 219 ---
$ws.Range("A219").Value = 45110.33333333334
$ws.Range("B219").Value = "FX_IDC:USDRON"
$ws.Range("C219").Value = 4.5517
$ws.Range("D219").Value = 4.5716
$ws.Range("E219").Value = 4.4849
$ws.Range("F219").Value = 4.4849
$ws.Range("G219").Value = 0
